$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new monthly rows (01-07-2021 and 01-08-2021) to the data table,
# mirroring the existing rows' layout: Serie (date label, text) + 9 numeric
# columns (B..J).
#
# Column A holds text labels such as "01-06-2021" that look like dates; a
# plain .Value assignment of such a literal is auto-interpreted by Excel as
# a date serial number. To keep these as plain text (shared strings), we
# briefly force a text number format before assigning the value and then
# clear the formatting again so the cell ends up with the default style,
# exactly like the rest of column A.

# Row 152: 01-07-2021
$ws.Range("A152").NumberFormat = "@"
$ws.Range("A152").Value = "01-07-2021"
$ws.Range("A152").ClearFormats()
$ws.Range("B152").Value = 11940
$ws.Range("C152").Value = 1513
$ws.Range("D152").Value = 2261
$ws.Range("E152").Value = 613
$ws.Range("F152").Value = 1509
$ws.Range("G152").Value = 715
$ws.Range("H152").Value = 2486
$ws.Range("I152").Value = 1599
$ws.Range("J152").Value = 1243

# Row 153: 01-08-2021
$ws.Range("A153").NumberFormat = "@"
$ws.Range("A153").Value = "01-08-2021"
$ws.Range("A153").ClearFormats()
$ws.Range("B153").Value = 16535
$ws.Range("C153").Value = 1828
$ws.Range("D153").Value = 2508
$ws.Range("E153").Value = 1904
$ws.Range("F153").Value = 3334
$ws.Range("G153").Value = 1485
$ws.Range("H153").Value = 3335
$ws.Range("I153").Value = 1264
$ws.Range("J153").Value = 878
